$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2)
$run = $para.Runs(2)

$run.Text = "My friend Brian, an engineer at One Medical, is working on an app that helps medical practices automate these types of tasks. He will provide a dataset containing few hundred labeled medical documents (as images) along with their text. We will select a subset of these types of documents to build a classifier (such as medical history records, lab results and release forms). Each "

$run2 = $run.InsertAfter("document page ")
$run3 = $run2.InsertAfter("is the unit of observation, and features we will extract will contain parsed text and/or visual aspects of the images. ")
